# New Changes To xlsx
# Rebuilds Sheet1 with the Manager_id / emp_dept / emp_share(%) table
# (mirrored across columns A:C and D:F), matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (shared string) even when it looks like a
# number (e.g. "1001", " 40") without leaving a lingering NumberFormat/style
# on the cell - match real Excel's plain "typed as text" cells.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 1 - headers (duplicated A:C -> D:F)
Set-TextValue $ws.Range("A1") "Manager_id  "
Set-TextValue $ws.Range("B1") "emp_dept"
Set-TextValue $ws.Range("C1") "emp_share (%)"
Set-TextValue $ws.Range("D1") "Manager_id  "
Set-TextValue $ws.Range("E1") "emp_dept"
Set-TextValue $ws.Range("F1") "emp_share (%)"

# Row 2
Set-TextValue $ws.Range("A2") "Null"
$ws.Range("B2").Value = "Finance"
Set-TextValue $ws.Range("C2") "60"
Set-TextValue $ws.Range("D2") "Null"
$ws.Range("E2").Value = "Finance"
Set-TextValue $ws.Range("F2") "60"

# Row 3
Set-TextValue $ws.Range("A3") "1001"
$ws.Range("B3").Value = "Finance"
Set-TextValue $ws.Range("C3") "20"
Set-TextValue $ws.Range("D3") "1001"
$ws.Range("E3").Value = "Finance"
Set-TextValue $ws.Range("F3") "20"

# Row 4
Set-TextValue $ws.Range("A4") "1004"
$ws.Range("B4").Value = "R&D"
Set-TextValue $ws.Range("C4") "30"
Set-TextValue $ws.Range("D4") "1004"
$ws.Range("E4").Value = "R&D"
Set-TextValue $ws.Range("F4") "30"

# Row 5
Set-TextValue $ws.Range("A5") "1004"
$ws.Range("B5").Value = "R&D"
Set-TextValue $ws.Range("C5") " 40"
Set-TextValue $ws.Range("D5") "1004"
$ws.Range("E5").Value = "R&D"
Set-TextValue $ws.Range("F5") " 40"

# Row 6
Set-TextValue $ws.Range("A6") "1001"
$ws.Range("B6").Value = "Finance"
Set-TextValue $ws.Range("C6") " 20"
Set-TextValue $ws.Range("D6") "1001"
$ws.Range("E6").Value = "Finance"
Set-TextValue $ws.Range("F6") " 20"

# Row 7
Set-TextValue $ws.Range("A7") "1005"
$ws.Range("B7").Value = "Finance"
Set-TextValue $ws.Range("C7") "15"
Set-TextValue $ws.Range("D7") "1005"
$ws.Range("E7").Value = "Finance"
Set-TextValue $ws.Range("F7") "15"

# Row 8
Set-TextValue $ws.Range("A8") "1001"
$ws.Range("B8").Value = "Finance"
Set-TextValue $ws.Range("C8") "25"
Set-TextValue $ws.Range("D8") "1001"
$ws.Range("E8").Value = "Finance"
Set-TextValue $ws.Range("F8") "25"

# Match the saved selection from the target file
$ws.Range("C7").Select()
